# Auto-generated script applying the BRVM automatic data refresh
# described by the commit "🔄 MAJ automatique BRVM via GitHub Actions".
$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations" (sector + stock recommendation table) ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 1: Titre
$ws1.Cells.Item(1, 1).Value = "Titre"
$ws1.Cells.Item(1, 2).Value = "Jours en Hausse"
$ws1.Cells.Item(1, 3).Value = "Jours en Baisse"
$ws1.Cells.Item(1, 4).Value = "Variation Totale (%)"
$ws1.Cells.Item(1, 5).Value = "Dernière Variation (%)"
$ws1.Cells.Item(1, 6).Value = "Recommandation"
$ws1.Cells.Item(1, 7).Value = "Stratégie"

# Row 2: BRVM - SERVICES PUBLICS
$ws1.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 7
$ws1.Cells.Item(2, 4).Value = 2564.27
$ws1.Cells.Item(2, 5).Value = 107.19
$ws1.Cells.Item(2, 6).Value = "🟡 Observer"
$ws1.Cells.Item(2, 7).Value = "➖ Neutre"

# Row 3: BRVM - AUTRES SECTEURS
$ws1.Cells.Item(3, 1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 3
$ws1.Cells.Item(3, 4).Value = 2070.4
$ws1.Cells.Item(3, 5).Value = 669.2
$ws1.Cells.Item(3, 6).Value = "🟡 Observer"
$ws1.Cells.Item(3, 7).Value = "➖ Neutre"

# Row 4: BRVM - DISTRIBUTION
$ws1.Cells.Item(4, 1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 3
$ws1.Cells.Item(4, 4).Value = 1490
$ws1.Cells.Item(4, 5).Value = 494.78
$ws1.Cells.Item(4, 6).Value = "🟡 Observer"
$ws1.Cells.Item(4, 7).Value = "➖ Neutre"

# Row 5: BRVM - TRANSPORT
$ws1.Cells.Item(5, 1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 3
$ws1.Cells.Item(5, 4).Value = 1097.81
$ws1.Cells.Item(5, 5).Value = 370.83
$ws1.Cells.Item(5, 6).Value = "🟡 Observer"
$ws1.Cells.Item(5, 7).Value = "➖ Neutre"

# Row 6: BRVM - AGRICULTURE
$ws1.Cells.Item(6, 1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 3
$ws1.Cells.Item(6, 4).Value = 1005.8
$ws1.Cells.Item(6, 5).Value = 341.63
$ws1.Cells.Item(6, 6).Value = "🟡 Observer"
$ws1.Cells.Item(6, 7).Value = "➖ Neutre"

# Row 7: BRVM - INDUSTRIE  (**)
$ws1.Cells.Item(7, 1).Value = "BRVM - INDUSTRIE  (**)"
$ws1.Cells.Item(7, 2).Value = 0
$ws1.Cells.Item(7, 3).Value = 3
$ws1.Cells.Item(7, 4).Value = 767.31
$ws1.Cells.Item(7, 5).Value = 257.22
$ws1.Cells.Item(7, 6).Value = "🟡 Observer"
$ws1.Cells.Item(7, 7).Value = "➖ Neutre"

# Row 8: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Cells.Item(8, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(8, 2).Value = 0
$ws1.Cells.Item(8, 3).Value = 4
$ws1.Cells.Item(8, 4).Value = 694.96
$ws1.Cells.Item(8, 5).Value = 170.29
$ws1.Cells.Item(8, 6).Value = "🟡 Observer"
$ws1.Cells.Item(8, 7).Value = "➖ Neutre"

# Row 9: BRVM-PRINCIPAL  (**)
$ws1.Cells.Item(9, 1).Value = "BRVM-PRINCIPAL  (**)"
$ws1.Cells.Item(9, 2).Value = 0
$ws1.Cells.Item(9, 3).Value = 3
$ws1.Cells.Item(9, 4).Value = 649.99
$ws1.Cells.Item(9, 5).Value = 217.65
$ws1.Cells.Item(9, 6).Value = "🟡 Observer"
$ws1.Cells.Item(9, 7).Value = "➖ Neutre"

# Row 10: BRVM - CONSOMMATION DE BASE  (**)
$ws1.Cells.Item(10, 1).Value = "BRVM - CONSOMMATION DE BASE  (**)"
$ws1.Cells.Item(10, 2).Value = 0
$ws1.Cells.Item(10, 3).Value = 3
$ws1.Cells.Item(10, 4).Value = 641.94
$ws1.Cells.Item(10, 5).Value = 216.08
$ws1.Cells.Item(10, 6).Value = "🟡 Observer"
$ws1.Cells.Item(10, 7).Value = "➖ Neutre"

# Row 11: BRVM - SERVICES FINANCIERS
$ws1.Cells.Item(11, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(11, 2).Value = 0
$ws1.Cells.Item(11, 3).Value = 4
$ws1.Cells.Item(11, 4).Value = 580.38
$ws1.Cells.Item(11, 5).Value = 145.91
$ws1.Cells.Item(11, 6).Value = "🟡 Observer"
$ws1.Cells.Item(11, 7).Value = "➖ Neutre"

# Row 12: BRVM-PRESTIGE
$ws1.Cells.Item(12, 1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(12, 2).Value = 0
$ws1.Cells.Item(12, 3).Value = 4
$ws1.Cells.Item(12, 4).Value = 570.66
$ws1.Cells.Item(12, 5).Value = 143.48
$ws1.Cells.Item(12, 6).Value = "🟡 Observer"
$ws1.Cells.Item(12, 7).Value = "➖ Neutre"

# Row 13: BRVM - INDUSTRIELS
$ws1.Cells.Item(13, 1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(13, 2).Value = 0
$ws1.Cells.Item(13, 3).Value = 4
$ws1.Cells.Item(13, 4).Value = 529.34
$ws1.Cells.Item(13, 5).Value = 132.35
$ws1.Cells.Item(13, 6).Value = "🟡 Observer"
$ws1.Cells.Item(13, 7).Value = "➖ Neutre"

# Row 14: BRVM - ENERGIE
$ws1.Cells.Item(14, 1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(14, 2).Value = 0
$ws1.Cells.Item(14, 3).Value = 4
$ws1.Cells.Item(14, 4).Value = 443.25
$ws1.Cells.Item(14, 5).Value = 110.99
$ws1.Cells.Item(14, 6).Value = "🟡 Observer"
$ws1.Cells.Item(14, 7).Value = "➖ Neutre"

# Row 15: BRVM - FINANCES
$ws1.Cells.Item(15, 1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(15, 2).Value = 0
$ws1.Cells.Item(15, 3).Value = 3
$ws1.Cells.Item(15, 4).Value = 442.08
$ws1.Cells.Item(15, 5).Value = 148.41
$ws1.Cells.Item(15, 6).Value = "🟡 Observer"
$ws1.Cells.Item(15, 7).Value = "➖ Neutre"

# Row 16: BRVM - TELECOMMUNICATIONS
$ws1.Cells.Item(16, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(16, 2).Value = 0
$ws1.Cells.Item(16, 3).Value = 4
$ws1.Cells.Item(16, 4).Value = 373.71
$ws1.Cells.Item(16, 5).Value = 94.05
$ws1.Cells.Item(16, 6).Value = "🟡 Observer"
$ws1.Cells.Item(16, 7).Value = "➖ Neutre"

# Row 17: BRVM-PRINCIPAL    (**)
$ws1.Cells.Item(17, 1).Value = "BRVM-PRINCIPAL    (**)"
$ws1.Cells.Item(17, 2).Value = 0
$ws1.Cells.Item(17, 3).Value = 1
$ws1.Cells.Item(17, 4).Value = 217.75
$ws1.Cells.Item(17, 5).Value = 217.75
$ws1.Cells.Item(17, 6).Value = "🟡 Observer"
$ws1.Cells.Item(17, 7).Value = "➖ Neutre"

# Row 18: BRVM - CONSOMMATION DE BASE    (**)
$ws1.Cells.Item(18, 1).Value = "BRVM - CONSOMMATION DE BASE    (**)"
$ws1.Cells.Item(18, 2).Value = 0
$ws1.Cells.Item(18, 3).Value = 1
$ws1.Cells.Item(18, 4).Value = 215.92
$ws1.Cells.Item(18, 5).Value = 215.92
$ws1.Cells.Item(18, 6).Value = "🟡 Observer"
$ws1.Cells.Item(18, 7).Value = "➖ Neutre"

# Row 19: BRVM – COMPOSITE TOTAL RETURN    (**)
$ws1.Cells.Item(19, 1).Value = "BRVM – COMPOSITE TOTAL RETURN    (**)"
$ws1.Cells.Item(19, 2).Value = 0
$ws1.Cells.Item(19, 3).Value = 1
$ws1.Cells.Item(19, 4).Value = 132.65
$ws1.Cells.Item(19, 5).Value = 132.65
$ws1.Cells.Item(19, 6).Value = "🟡 Observer"
$ws1.Cells.Item(19, 7).Value = "➖ Neutre"

# Row 20: EVIOSYS PACKAGING SIEM CI (SEMC)
$ws1.Cells.Item(20, 1).Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Cells.Item(20, 2).Value = 1
$ws1.Cells.Item(20, 3).Value = 0
$ws1.Cells.Item(20, 4).Value = 7.14
$ws1.Cells.Item(20, 5).Value = 7.14
$ws1.Cells.Item(20, 6).Value = "🟡 Observer"
$ws1.Cells.Item(20, 7).Value = "➖ Neutre"

# Row 21: SONATEL SN (SNTS)
$ws1.Cells.Item(21, 1).Value = "SONATEL SN (SNTS)"
$ws1.Cells.Item(21, 2).Value = 1
$ws1.Cells.Item(21, 3).Value = 0
$ws1.Cells.Item(21, 4).Value = 7.05
$ws1.Cells.Item(21, 5).Value = 7.05
$ws1.Cells.Item(21, 6).Value = "🟡 Observer"
$ws1.Cells.Item(21, 7).Value = "➖ Neutre"

# Row 22: NEI-CEDA CI (NEIC)
$ws1.Cells.Item(22, 1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(22, 2).Value = 2
$ws1.Cells.Item(22, 3).Value = 1
$ws1.Cells.Item(22, 4).Value = 6.44
$ws1.Cells.Item(22, 5).Value = -7.5
$ws1.Cells.Item(22, 6).Value = "🟡 Observer"
$ws1.Cells.Item(22, 7).Value = "👀 À surveiller"

# Row 23: SAPH CI (SPHC)
$ws1.Cells.Item(23, 1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(23, 2).Value = 1
$ws1.Cells.Item(23, 3).Value = 0
$ws1.Cells.Item(23, 4).Value = 5.33
$ws1.Cells.Item(23, 5).Value = 5.33
$ws1.Cells.Item(23, 6).Value = "🟡 Observer"
$ws1.Cells.Item(23, 7).Value = "➖ Neutre"

# Row 24: SOCIETE GENERALE COTE D'IVOIRE (SGBC)
$ws1.Cells.Item(24, 1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(24, 2).Value = 1
$ws1.Cells.Item(24, 3).Value = 0
$ws1.Cells.Item(24, 4).Value = 4.75
$ws1.Cells.Item(24, 5).Value = 4.75
$ws1.Cells.Item(24, 6).Value = "🟡 Observer"
$ws1.Cells.Item(24, 7).Value = "➖ Neutre"

# Row 25: CIE CI (CIEC)
$ws1.Cells.Item(25, 1).Value = "CIE CI (CIEC)"
$ws1.Cells.Item(25, 2).Value = 1
$ws1.Cells.Item(25, 3).Value = 0
$ws1.Cells.Item(25, 4).Value = 4.66
$ws1.Cells.Item(25, 5).Value = 4.66
$ws1.Cells.Item(25, 6).Value = "🟡 Observer"
$ws1.Cells.Item(25, 7).Value = "➖ Neutre"

# Row 26: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Cells.Item(26, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(26, 2).Value = 1
$ws1.Cells.Item(26, 3).Value = 0
$ws1.Cells.Item(26, 4).Value = 4.55
$ws1.Cells.Item(26, 5).Value = 4.55
$ws1.Cells.Item(26, 6).Value = "🟡 Observer"
$ws1.Cells.Item(26, 7).Value = "➖ Neutre"

# Row 27: VIVO ENERGY CI (SHEC)
$ws1.Cells.Item(27, 1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(27, 2).Value = 1
$ws1.Cells.Item(27, 3).Value = 0
$ws1.Cells.Item(27, 4).Value = 3.7
$ws1.Cells.Item(27, 5).Value = 3.7
$ws1.Cells.Item(27, 6).Value = "🟡 Observer"
$ws1.Cells.Item(27, 7).Value = "➖ Neutre"

# Row 28: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws1.Cells.Item(28, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(28, 2).Value = 1
$ws1.Cells.Item(28, 3).Value = 0
$ws1.Cells.Item(28, 4).Value = 3.42
$ws1.Cells.Item(28, 5).Value = 3.42
$ws1.Cells.Item(28, 6).Value = "🟡 Observer"
$ws1.Cells.Item(28, 7).Value = "➖ Neutre"

# Row 29: SOCIETE IVOIRIENNE DE BANQUE  (SIBC)
$ws1.Cells.Item(29, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Cells.Item(29, 2).Value = 1
$ws1.Cells.Item(29, 3).Value = 0
$ws1.Cells.Item(29, 4).Value = 3.23
$ws1.Cells.Item(29, 5).Value = 3.23
$ws1.Cells.Item(29, 6).Value = "🟡 Observer"
$ws1.Cells.Item(29, 7).Value = "➖ Neutre"

# Row 30: SICOR CI (SICC)
$ws1.Cells.Item(30, 1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(30, 2).Value = 1
$ws1.Cells.Item(30, 3).Value = 1
$ws1.Cells.Item(30, 4).Value = 1
$ws1.Cells.Item(30, 5).Value = -5.71
$ws1.Cells.Item(30, 6).Value = "🟡 Observer"
$ws1.Cells.Item(30, 7).Value = "👀 À surveiller"

# Row 31: TRACTAFRIC MOTORS CI (PRSC)
$ws1.Cells.Item(31, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(31, 2).Value = 1
$ws1.Cells.Item(31, 3).Value = 1
$ws1.Cells.Item(31, 4).Value = 0.81
$ws1.Cells.Item(31, 5).Value = -3.43
$ws1.Cells.Item(31, 6).Value = "🟡 Observer"
$ws1.Cells.Item(31, 7).Value = "👀 À surveiller"

# Row 32: NESTLE CI (NTLC)
$ws1.Cells.Item(32, 1).Value = "NESTLE CI (NTLC)"
$ws1.Cells.Item(32, 2).Value = 1
$ws1.Cells.Item(32, 3).Value = 1
$ws1.Cells.Item(32, 4).Value = 0.46
$ws1.Cells.Item(32, 5).Value = 3.24
$ws1.Cells.Item(32, 6).Value = "🟡 Observer"
$ws1.Cells.Item(32, 7).Value = "👀 À surveiller"

# Row 33: TOTAL
$ws1.Cells.Item(33, 1).Value = "TOTAL"
$ws1.Cells.Item(33, 2).Value = 0
$ws1.Cells.Item(33, 3).Value = 2
$ws1.Cells.Item(33, 4).Value = 0
$ws1.Cells.Item(33, 5).Value = 0
$ws1.Cells.Item(33, 6).Value = "🟡 Observer"
$ws1.Cells.Item(33, 7).Value = "➖ Neutre"

# Row 34: ORAGROUP TOGO (ORGT)
$ws1.Cells.Item(34, 1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(34, 2).Value = 2
$ws1.Cells.Item(34, 3).Value = 2
$ws1.Cells.Item(34, 4).Value = -0.31
$ws1.Cells.Item(34, 5).Value = 3.54
$ws1.Cells.Item(34, 6).Value = "🟡 Observer"
$ws1.Cells.Item(34, 7).Value = "👀 À surveiller"

# Row 35: BANK OF AFRICA BF (BOABF)
$ws1.Cells.Item(35, 1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(35, 2).Value = 1
$ws1.Cells.Item(35, 3).Value = 1
$ws1.Cells.Item(35, 4).Value = -0.95
$ws1.Cells.Item(35, 5).Value = 5.07
$ws1.Cells.Item(35, 6).Value = "🟡 Observer"
$ws1.Cells.Item(35, 7).Value = "👀 À surveiller"

# Row 36: ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)
$ws1.Cells.Item(36, 1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Cells.Item(36, 2).Value = 2
$ws1.Cells.Item(36, 3).Value = 2
$ws1.Cells.Item(36, 4).Value = -1.57
$ws1.Cells.Item(36, 5).Value = 6.33
$ws1.Cells.Item(36, 6).Value = "🟡 Observer"
$ws1.Cells.Item(36, 7).Value = "👀 À surveiller"

# Row 37: BERNABE CI (BNBC)
$ws1.Cells.Item(37, 1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(37, 2).Value = 0
$ws1.Cells.Item(37, 3).Value = 1
$ws1.Cells.Item(37, 4).Value = -3.19
$ws1.Cells.Item(37, 5).Value = -3.19
$ws1.Cells.Item(37, 6).Value = "🟡 Observer"
$ws1.Cells.Item(37, 7).Value = "➖ Neutre"

# Row 38: BANK OF AFRICA ML (BOAM)
$ws1.Cells.Item(38, 1).Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(38, 2).Value = 0
$ws1.Cells.Item(38, 3).Value = 1
$ws1.Cells.Item(38, 4).Value = -3.5
$ws1.Cells.Item(38, 5).Value = -3.5
$ws1.Cells.Item(38, 6).Value = "🟡 Observer"
$ws1.Cells.Item(38, 7).Value = "➖ Neutre"

# Row 39: LOTERIE NATIONALE DU BENIN (LNBB)
$ws1.Cells.Item(39, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(39, 2).Value = 0
$ws1.Cells.Item(39, 3).Value = 1
$ws1.Cells.Item(39, 4).Value = -3.65
$ws1.Cells.Item(39, 5).Value = -3.65
$ws1.Cells.Item(39, 6).Value = "🟡 Observer"
$ws1.Cells.Item(39, 7).Value = "➖ Neutre"

# Row 40: CFAO MOTORS CI (CFAC)
$ws1.Cells.Item(40, 1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(40, 2).Value = 0
$ws1.Cells.Item(40, 3).Value = 2
$ws1.Cells.Item(40, 4).Value = -3.83
$ws1.Cells.Item(40, 5).Value = -2.45
$ws1.Cells.Item(40, 6).Value = "🟡 Observer"
$ws1.Cells.Item(40, 7).Value = "➖ Neutre"

# Row 41: SAFCA CI (SAFC)
$ws1.Cells.Item(41, 1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(41, 2).Value = 0
$ws1.Cells.Item(41, 3).Value = 1
$ws1.Cells.Item(41, 4).Value = -3.93
$ws1.Cells.Item(41, 5).Value = -3.93
$ws1.Cells.Item(41, 6).Value = "🟡 Observer"
$ws1.Cells.Item(41, 7).Value = "➖ Neutre"

# Row 42: SERVAIR ABIDJAN CI (ABJC)
$ws1.Cells.Item(42, 1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(42, 2).Value = 0
$ws1.Cells.Item(42, 3).Value = 1
$ws1.Cells.Item(42, 4).Value = -4.92
$ws1.Cells.Item(42, 5).Value = -4.92
$ws1.Cells.Item(42, 6).Value = "🟡 Observer"
$ws1.Cells.Item(42, 7).Value = "➖ Neutre"

# Row 43: SUCRIVOIRE (SCRC)
$ws1.Cells.Item(43, 1).Value = "SUCRIVOIRE (SCRC)"
$ws1.Cells.Item(43, 2).Value = 0
$ws1.Cells.Item(43, 3).Value = 1
$ws1.Cells.Item(43, 4).Value = -5.45
$ws1.Cells.Item(43, 5).Value = -5.45
$ws1.Cells.Item(43, 6).Value = "🟡 Observer"
$ws1.Cells.Item(43, 7).Value = "➖ Neutre"

# Row 44: SETAO CI (STAC)
$ws1.Cells.Item(44, 1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(44, 2).Value = 1
$ws1.Cells.Item(44, 3).Value = 3
$ws1.Cells.Item(44, 4).Value = -9.91
$ws1.Cells.Item(44, 5).Value = -3.4
$ws1.Cells.Item(44, 6).Value = "🔴 Vente"
$ws1.Cells.Item(44, 7).Value = "⚠️ Risque de décrochage"

# --- Sheet "Top_YTD" (top sector performers by YTD progression) ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 1: Titre
$ws2.Cells.Item(1, 1).Value = "Titre"
$ws2.Cells.Item(1, 2).Value = "Progression YTD (%)"

# Row 2: BRVM-PRINCIPAL    (**)
$ws2.Cells.Item(2, 1).Value = "BRVM-PRINCIPAL    (**)"
$ws2.Cells.Item(2, 2).Value = 217.75

# Row 3: BRVM - CONSOMMATION DE BASE    (**)
$ws2.Cells.Item(3, 1).Value = "BRVM - CONSOMMATION DE BASE    (**)"
$ws2.Cells.Item(3, 2).Value = 215.92

# Row 4: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws2.Cells.Item(4, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Cells.Item(4, 2).Value = 170.29

# Row 5: BRVM - SERVICES FINANCIERS
$ws2.Cells.Item(5, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws2.Cells.Item(5, 2).Value = 145.91

# Row 6: BRVM-PRESTIGE
$ws2.Cells.Item(6, 1).Value = "BRVM-PRESTIGE"
$ws2.Cells.Item(6, 2).Value = 143.48

# Row 7: BRVM – COMPOSITE TOTAL RETURN    (**)
$ws2.Cells.Item(7, 1).Value = "BRVM – COMPOSITE TOTAL RETURN    (**)"
$ws2.Cells.Item(7, 2).Value = 132.65

# Row 8: BRVM - INDUSTRIELS
$ws2.Cells.Item(8, 1).Value = "BRVM - INDUSTRIELS"
$ws2.Cells.Item(8, 2).Value = 132.35

# Row 9: BRVM - ENERGIE
$ws2.Cells.Item(9, 1).Value = "BRVM - ENERGIE"
$ws2.Cells.Item(9, 2).Value = 110.99

# Row 10: BRVM - SERVICES PUBLICS
$ws2.Cells.Item(10, 1).Value = "BRVM - SERVICES PUBLICS"
$ws2.Cells.Item(10, 2).Value = 107.19

# Row 11: BRVM - TELECOMMUNICATIONS
$ws2.Cells.Item(11, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Cells.Item(11, 2).Value = 94.05

